# Rename header row suffixes:
#   "<Name>_old" -> "<Name>_FV2410"
#   "<Name>_new" -> "<Name>_FV2504"
# then wrap the data range in an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID",
    "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung"
)

# Columns A-J (1-10): "<Name>_old" -> "<Name>_FV2410"
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i] + "_FV2410"
}

# Column K (11) "diff" stays the same.

# Columns L-U (12-21): "<Name>_new" -> "<Name>_FV2504"
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 12)
    $cell.Value = $headers[$i] + "_FV2504"
}

# Wrap the whole data range into an Excel Table (adds xl/tables/table1.xml
# and the <tableParts> reference on the worksheet).
$dataRange = $ws.Range("A1:U84")
$lo = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$lo.Name = "Table1"

# Freeze the header row (pane split after row 1).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
